$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("F2").Value = 3
$ws.Range("F8").Value = -3
$ws.Range("F18").Value = -1
$ws.Range("F36").Value = 1
$ws.Range("F40").Value = 0
$ws.Range("F46").Value = 4
$ws.Range("F52").Value = -1
$ws.Range("F60").Value = -4
$ws.Range("F67").Value = -9
